$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.005.90"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -2.23%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.816.21"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.44%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -1.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.50"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.61%  "

$ws.Range("E6").Value = "  -1.07%  "

$ws.Range("E7").Value = "  -2.29%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3666"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.85%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07201"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.81%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8406"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.84%  "

$ws.Range("E11").Value = "  -3.58%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.818.07"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.16%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.623"
$ws.Range("D13").ClearFormats()

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07078"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.70%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.271"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -3.03%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.74"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.55%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.26%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008787"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.96%  "

$ws.Range("E19").Value = "  -1.10%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.94"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -3.11%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.078.00"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.99%  "

$ws.Range("E22").Value = "  -1.86%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.79"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.51%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.042.21"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.88%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.975"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.67%  "

$ws.Range("E26").Value = "  -2.62%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.218"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +3.11%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.26"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.70%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.196"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.99%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.77"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.24%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08783"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.72%  "

$ws.Range("E32").Value = "  -4.33%  "

$ws.Range("E33").Value = "  +2.38%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7369"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -4.78%  "

$ws.Range("E35").Value = "  -2.93%  "

$ws.Range("E36").Value = "  -1.16%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.092"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.53%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01953"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.87%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05226"
$ws.Range("D39").ClearFormats()

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.262"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.12%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.865"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.94%  "

$ws.Range("E42").Value = "  +0.11%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5025"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.52%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.556"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.50%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.56"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.81%  "

$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "106.08"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.72%  "

$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4727"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.11%  "

$ws.Range("E48").Value = "  -1.19%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06354"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.90%  "

$ws.Range("E50").Value = "  -2.81%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.872"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.71%  "
